$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell C3: "№ студенческого билета, зачетки" -> "N_ZACHET"
$ws.Range("C3").Value = "N_ZACHET"

# Move the active selection to D7, matching the saved view state in the diff
$ws.Range("D7").Select()
